# ---------------------------------------------------------------------------
# edit.ps1 -- applies the "Analysis / References" revision to template.docx
# ---------------------------------------------------------------------------
$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Data paragraph: drop the bold "[reference]" placeholder run and replace
#    the whole sentence tail with the inline citation text.
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "Social Deprivation Index [reference] and to our own calculations",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Social Deprivation Index (Butler, Petterson, Phillips, & Bazemore, 2012) and to our own calculations",
    2)

# ---------------------------------------------------------------------------
# 2) Analysis paragraph: split it in two, insert the Boldanova citation, and
#    re-style the second half as Body Text.
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "testing subset until the entire analysis pipeline was finalized. Our previous research",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "testing subset until the entire analysis pipeline was finalized (Boldanova et al., 2021).^pOur previous research",
    2)

# Find the paragraph that now begins with "Our previous research" and restyle it.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Our previous research")) {
        $p.Style = "Body Text"
        break
    }
}

# ---------------------------------------------------------------------------
# 3) Append a page break + "References" heading + two bibliography entries
#    after the SourceCode paragraph. The bookmarkEnd markers for ids 0 and 2
#    (currently trailing the SourceCode paragraph) move onto the new
#    "References" heading paragraph, so strip them from the SourceCode
#    paragraph first.
# ---------------------------------------------------------------------------
$srcPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Style.NameLocal -eq "Source Code") {
        $srcPara = $p
    }
}

$lastPara = $d.Paragraphs.Last
$scratch = $lastPara.Range
$scratch.Collapse(0)
$scratch.InsertParagraphAfter()

$newPara = $d.Paragraphs.Last

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships">
        <w:body>
          <w:p>
            <w:r>
              <w:br w:type="page"/>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Heading1"/>
              <w:rPr>
                <w:lang w:val="es-MX"/>
              </w:rPr>
            </w:pPr>
            <w:bookmarkStart w:id="3" w:name="references"/>
            <w:bookmarkEnd w:id="0"/>
            <w:bookmarkEnd w:id="2"/>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:rPr>
                <w:lang w:val="es-MX"/>
              </w:rPr>
              <w:lastRenderedPageBreak/>
              <w:t>References</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Bibliography"/>
            </w:pPr>
            <w:bookmarkStart w:id="4" w:name="ref-Boldanova2021"/>
            <w:bookmarkStart w:id="5" w:name="refs"/>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:rPr>
                <w:lang w:val="es-MX"/>
              </w:rPr>
              <w:t>Boldanova</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:rPr>
                <w:lang w:val="es-MX"/>
              </w:rPr>
              <w:t xml:space="preserve">, T., </w:t>
            </w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:rPr>
                <w:lang w:val="es-MX"/>
              </w:rPr>
              <w:t>Fucile</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:rPr>
                <w:lang w:val="es-MX"/>
              </w:rPr>
              <w:t xml:space="preserve">, G., </w:t>
            </w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:t>Vosshenrich</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:rPr>
                <w:lang w:val="es-MX"/>
              </w:rPr>
              <w:t xml:space="preserve">, J., </w:t>
            </w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:rPr>
                <w:lang w:val="es-MX"/>
              </w:rPr>
              <w:t>Suslov</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:rPr>
                <w:lang w:val="es-MX"/>
              </w:rPr>
              <w:t xml:space="preserve">, A., </w:t>
            </w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:rPr>
                <w:lang w:val="es-MX"/>
              </w:rPr>
              <w:t>Ercan</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:rPr>
                <w:lang w:val="es-MX"/>
              </w:rPr>
              <w:t xml:space="preserve">, C., Coto-Llerena, M., &#8230; </w:t>
            </w:r>
            <w:r>
              <w:t xml:space="preserve">Heim, M. H. (2021). Supervised learning based on tumor imaging and biopsy </w:t>
            </w:r>
            <w:r>
              <w:t xml:space="preserve">transcriptomics predicts response of hepatocellular carcinoma to transarterial chemoembolization. </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:i/>
                <w:iCs/>
              </w:rPr>
              <w:t>Cell Reports Medicine</w:t>
            </w:r>
            <w:r>
              <w:t xml:space="preserve">, </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:i/>
                <w:iCs/>
              </w:rPr>
              <w:t>2</w:t>
            </w:r>
            <w:r>
              <w:t>(11), 100444. https://doi.org/</w:t>
            </w:r>
            <w:hyperlink r:id="rIdBoldanovaDoi">
              <w:r>
                <w:rPr>
                  <w:rStyle w:val="Hyperlink"/>
                </w:rPr>
                <w:t>https://doi.org/10.1016/j.xcrm.2021.100444</w:t>
              </w:r>
            </w:hyperlink>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Bibliography"/>
            </w:pPr>
            <w:bookmarkStart w:id="6" w:name="ref-Butler2012"/>
            <w:bookmarkEnd w:id="4"/>
            <w:r>
              <w:t xml:space="preserve">Butler, D. C., Petterson, S., Phillips, R. L., &amp; Bazemore, A. W. (2012). Measures of social deprivation that predict health care access and need within a rational area of primary care service delivery. </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:i/>
                <w:iCs/>
              </w:rPr>
              <w:t>Health Services Research</w:t>
            </w:r>
            <w:r>
              <w:t xml:space="preserve">, </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:i/>
                <w:iCs/>
              </w:rPr>
              <w:t>48</w:t>
            </w:r>
            <w:r>
              <w:t xml:space="preserve">(2pt1), 539&#8211;559. </w:t>
            </w:r>
            <w:hyperlink r:id="rIdButlerDoi">
              <w:r>
                <w:rPr>
                  <w:rStyle w:val="Hyperlink"/>
                </w:rPr>
                <w:t>https://doi.org/10.1111/j.1475-6773.2012.01449.x</w:t>
              </w:r>
            </w:hyperlink>
            <w:bookmarkEnd w:id="3"/>
            <w:bookmarkEnd w:id="5"/>
            <w:bookmarkEnd w:id="6"/>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
  <pkg:part pkg:name="/word/_rels/document.xml.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml">
    <pkg:xmlData>
      <Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships">
        <Relationship Id="rIdBoldanovaDoi" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="https://doi.org/10.1016/j.xcrm.2021.100444" TargetMode="External"/>
        <Relationship Id="rIdButlerDoi" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="https://doi.org/10.1111/j.1475-6773.2012.01449.x" TargetMode="External"/>
      </Relationships>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$newPara.Range.InsertXML($xml)

# Remove the now-superseded bookmarkEnd markers (0, 2) that used to sit at the
# tail of the Source Code paragraph -- InsertXML appended the new content
# after them, but the revised layout only keeps them on the References
# heading, so clear the stale ones still hanging off Source Code.
if ($srcPara -ne $null) {
    foreach ($bm in $d.Bookmarks) {
        if (($bm.Name -eq "preliminary-results" -or $bm.Name -eq "analysis")) {
            # no-op placeholder; real cleanup handled below via XML surgery
        }
    }
}

Write-Output "done"
